# Properties Listed and Hash registered
#
# Adds two new "hash" rows (A4, A5) below the existing three, using a small
# monospace font (Consolas 7pt, VS-Code-string orange FFCE9178) to set them
# apart visually from the rest of the list, then moves the active selection
# to K14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- append the two new hash values -------------------------------------
$ws.Range("A4").Value = "aabfebe72a7b7e802fc567d097979e1f728df9956b9e522016fef04903af903e"
$ws.Range("A5").Value = "19b0a0ce0dbfa2a6403592fffdf8c30587c6fbc78dae5eafa9ff803829b081c3"

# --- style the first new cell (Consolas 7pt, #CE9178) --------------------
$cell = $ws.Range("A4")
$cell.Font.Name = "Consolas"
$cell.Font.Size = 7
$cell.Font.Family = 3
$cell.Font.Color = 7901646

# --- carry the exact same formatting onto the second new cell ------------
$cell.Copy()
$ws.Range("A5").PasteSpecial(-4122)

# --- move the selection, matching where the user clicked next ------------
[void]$ws.Range("K14").Select()
